$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data-driven login test values: replace the single test/1234 credentials
# with two rows of credentials (mcAngular/3d6g4f7j5g8k and
# tmtmoney/oy06ri94uw73).
$ws.Range("A2").Value = "mcAngular"
$ws.Range("B2").Value = "3d6g4f7j5g8k"

$ws.Range("A3").Value = "tmtmoney"
$ws.Range("B3").Value = "oy06ri94uw73"

# The new B2 text no longer looks numeric, so it shouldn't carry the
# quote-prefix formatting the old "1234" value needed. Copy the plain
# (non quote-prefixed) format from A2 onto A3/B3 so the new row matches
# the rest of the table instead of inheriting default formatting.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
